# Apply recipe/consumption-rate changes for reactor coolant and cores.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 (Coolant Powder recipe): duration 7.5 -> 6 (usage formula recalculates automatically)
$ws.Range("D12").Value = 6

# Row 14 (Cloudy Coolant recipe): duration 30 -> 24 (usage formula recalculates automatically)
$ws.Range("D14").Value = 24

# Row 16 (Reactor Coolant recipe): update the "(change to 5000)" note to the final "5000 (5m^3)" value
$ws.Range("J16").Value = "5000 (5m^3)"

# Remove the now-stale helper note in J17
$ws.Range("J17").ClearContents()

# Update the active selection to reflect where the user last clicked
$ws.Range("J10").Select()
